# Update the cryptos list with refreshed price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($r, $Coin, $Link, $Price, $Volume)
    if ($Coin -ne $null)   { $ws.Cells.Item($r, 2).Value = $Coin }
    if ($Link -ne $null)   { $ws.Cells.Item($r, 3).Value = $Link }
    if ($Price -ne $null)  {
        $c = $ws.Cells.Item($r, 4)
        if ($Price -match '^[0-9]+(\.[0-9]+)?$') {
            # Looks like a plain number (e.g. "31.72"); force it to stay text
            # (matching the source workbook's inline-string cells) without
            # leaving a visible quote-prefix style behind.
            $c.Value = "'" + $Price
            $c.Style = "Normal"
        } else {
            $c.Value = $Price
        }
    }
    if ($Volume -ne $null) { $ws.Cells.Item($r, 5).Value = $Volume }
}

Set-Row 2  $null $null "34.308.07" "  +12.35%  "
Set-Row 3  $null $null "1.817.42"  "  +7.80%  "
Set-Row 4  $null $null $null       "  +0.43%  "
Set-Row 5  $null $null "231.13"    "  +4.45%  "
Set-Row 6  $null $null $null       "  +4.94%  "
Set-Row 7  $null $null $null       "  +0.51%  "
Set-Row 8  $null $null "31.72"     "  +3.56%  "
Set-Row 9  $null $null "46.23"     "  +4.32%  "
Set-Row 12 $null $null $null       "  +3.40%  "
Set-Row 13 $null $null "2.079.97"  "  +7.77%  "
Set-Row 14 $null $null "1.816.39"  "  +7.74%  "
Set-Row 15 $null $null $null       "  +4.86%  "
Set-Row 16 $null $null "34.286.50" "  +12.16%  "
Set-Row 17 $null $null "10.37"     "  -2.81%  "
Set-Row 18 $null $null $null       "  +8.30%  "
Set-Row 19 $null $null "70.77"     "  +7.23%  "
Set-Row 20 $null $null "261.01"    "  +5.39%  "
Set-Row 21 $null $null "0.0₃0749"  "  +4.36%  "
Set-Row 22 $null $null $null       "  +0.82%  "
Set-Row 23 $null $null $null       "  +3.39%  "
Set-Row 24 $null $null "4.40"      "  +2.81%  "
Set-Row 25 $null $null "2.20"      "  -0.50%  "
Set-Row 26 $null $null "161.35"    "  +2.04%  "
Set-Row 27 $null $null "16.80"     "  +5.35%  "
Set-Row 28 $null $null $null       "  +4.94%  "
Set-Row 29 $null $null $null       "  +5.65%  "
Set-Row 30 $null $null $null       "  +0.41%  "
Set-Row 31 $null $null "3.85"      "  +10.12%  "

Set-Row 32 "Hedera"       "https://coinranking.com/coin/jad286TjB+hedera-hbar"       "0.0518" "  +3.56%  "
Set-Row 33 "PancakeSwap"  "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"  "1.22"   "  +7.17%  "

Set-Row 34 $null $null "3.59"      "  +8.74%  "
Set-Row 35 $null $null "1.591.22"  "  +5.72%  "
Set-Row 36 $null $null $null       "  +6.22%  "
Set-Row 37 $null $null $null       "  +2.50%  "
Set-Row 38 $null $null "85.99"     "  +8.77%  "

Set-Row 39 "ImmutableX" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"    "0.630"  "  +7.86%  "
Set-Row 40 "VeChain"    "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"   "0.0188" "  +5.20%  "

Set-Row 41 $null $null "2.79"      "  +1.44%  "
Set-Row 42 $null $null "2.36"      "  +1.86%  "
Set-Row 43 $null $null $null       "  +8.26%  "
Set-Row 44 $null $null $null       "  +6.69%  "
Set-Row 45 $null $null "0.0521"    "  +3.38%  "
Set-Row 46 $null $null $null       "  +5.56%  "
Set-Row 47 $null $null "1.977.10"  "  +8.23%  "
Set-Row 48 $null $null "53.90"     "  +3.00%  "
Set-Row 49 $null $null "5.74"      "  +5.85%  "
Set-Row 50 $null $null $null       "  +0.54%  "

Set-Row 51 "BabyDogeCoin" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge" "0.0₆0123" "  +8.41%  "
